# Refresh the Price (D) and Volume(1h) (E) columns for rows 2-51
# on the cryptos worksheet with the latest scraped snapshot values.
#
# Price/Volume cells are stored as plain text in this sheet (prices use
# '.' as a thousands separator, e.g. '41.621.49', and volumes carry
# padding + '%'). Assigning a numeric-looking string to .Value would let
# Excel auto-convert it to a real number, so for any new price that would
# parse as a plain number we force the Text format first and restore the
# original (General/Normal) cell style afterwards so no visible formatting
# changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.621.49'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '2.459.43'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  +0.13%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '318.82'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.61%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '90.98'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("E8").Value = '  +0.08%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.503'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.32%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0854'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.16%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '32.50'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").Value = '2.839.63'
$ws.Range("E13").Value = '  -0.45%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.81'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.15%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.41'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.06%  '
$ws.Range("D16").Value = '2.483.24'
$ws.Range("E16").Value = '  +0.93%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.783'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '41.524.21'
$ws.Range("E18").Value = '  -0.26%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.37'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.66%  '
$ws.Range("D20").Value = '0.0₃0935'
$ws.Range("E20").Value = '  -3.52%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '71.69'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("E22").Value = '  -3.01%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '237.48'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.68%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.74'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.43%  '
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  +0.07%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '24.56'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("E28").Value = '  -1.60%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.64'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.99%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '35.97'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.60%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '157.39'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.88%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.38'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  -0.26%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0757'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.86%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '16.82'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -4.71%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.115'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.86'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.90%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.81'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("E40").Value = '  -0.41%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.98'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("E42").Value = '  -7.11%  '
$ws.Range("D43").Value = '1.989.61'
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("E44").Value = '  -1.63%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '18.45'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.96%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.93'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.13%  '
$ws.Range("E47").Value = '  +4.21%  '
$ws.Range("D48").Value = '2.717.93'
$ws.Range("E48").Value = '  +0.45%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '75.28'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.85%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '96.51'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.77%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '66.39'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.69%  '
